$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "dnasr281@gmail.com, System"
$newVal = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldVal) {
        $cell.Value2 = $newVal
    }
}
